$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$explicacao = "O conteúdo do tema aborda a gestão da Pandemia, a regulação de atividades de linha de frente ou outras atividades da sociedade brasileira que foram afetadas pela COVID-19, como a assistência social em decorrência ao período da crise sanitária e afins."

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $jVal = $ws.Cells.Item($r, 10).Value2
    if ($dVal -eq "Pandemia da COVID-19" -and ($jVal -eq $null -or $jVal -eq "")) {
        $ws.Cells.Item($r, 10).Value = $explicacao
    }
}
